$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2:F25
$bf = New-Object 'object[,]' 24,5
$bf[0,0]=1.131363195441736; $bf[0,1]=0.1326042677527255; $bf[0,2]=0.1275294616269278; $bf[0,3]=0.1186993804594541; $bf[0,4]=1.560539851896451;
$bf[1,0]=1.054654112222835; $bf[1,1]=0.120284905689843; $bf[1,2]=0.1273487956012396; $bf[1,3]=0.1193505311781822; $bf[1,4]=1.564899490733488;
$bf[2,0]=1.007785563639914; $bf[2,1]=0.1126673836268992; $bf[2,2]=0.1272696348334676; $bf[2,3]=0.1197759119407662; $bf[2,4]=1.568351032597818;
$bf[3,0]=0.9887458202979644; $bf[3,1]=0.1095499019420316; $bf[3,2]=0.1272453982941961; $bf[3,3]=0.1199557005955256; $bf[3,4]=1.569952506591456;
$bf[4,0]=0.9855879260251754; $bf[4,1]=0.1090314486779107; $bf[4,2]=0.1272418593395557; $bf[4,3]=0.1199859438417414; $bf[4,4]=1.570230208211811;
$bf[5,0]=1.007528543881847; $bf[5,1]=0.1126253937541151; $bf[5,2]=0.1272692754446965; $bf[5,3]=0.1197783105311396; $bf[5,4]=1.568371841140362;
$bf[6,0]=1.10486677463382; $bf[6,1]=0.1283677145816; $bf[6,2]=0.1274605942275144; $bf[6,3]=0.1189185972995608; $bf[6,4]=1.561882344309538;
$bf[7,0]=1.297526890840004; $bf[7,1]=0.1588098103021309; $bf[7,2]=0.1280864623944211; $bf[7,3]=0.117435039435585; $bf[7,4]=1.555298731523109;
$bf[8,0]=1.440101607350186; $bf[8,1]=0.1809096999601252; $bf[8,2]=0.1286974027138115; $bf[8,3]=0.1164676359234793; $bf[8,4]=1.554201924670053;
$bf[9,0]=1.505174222380333; $bf[9,1]=0.1909048998321339; $bf[9,2]=0.1290077952422379; $bf[9,3]=0.1160539867443005; $bf[9,4]=1.554514328754593;
$bf[10,0]=1.529845116386298; $bf[10,1]=0.1946813389985778; $bf[10,2]=0.1291299704824453; $bf[10,3]=0.1159011362672722; $bf[10,4]=1.554749201419327;
$bf[11,0]=1.524530519342022; $bf[11,1]=0.1938683970670638; $bf[11,2]=0.1291034521847436; $bf[11,3]=0.1159338869843669; $bf[11,4]=1.554693434392377;
$bf[12,0]=1.507203333258531; $bf[12,1]=0.1912157615704189; $bf[12,2]=0.1290177539433017; $bf[12,3]=0.1160413357578687; $bf[12,4]=1.554531316159839;
$bf[13,0]=1.496593701804045; $bf[13,1]=0.1895898305164678; $bf[13,2]=0.1289658640997686; $bf[13,3]=0.1161076444197207; $bf[13,4]=1.554447192243686;
$bf[14,0]=1.4358531449613; $bf[14,1]=0.1802553059717127; $bf[14,2]=0.1286777684249358; $bf[14,3]=0.1164951997383152; $bf[14,4]=1.554197826276422;
$bf[15,0]=1.398644649235962; $bf[15,1]=0.1745138669374171; $bf[15,2]=0.1285093244356261; $bf[15,3]=0.1167397134214543; $bf[15,4]=1.554252583662034;
$bf[16,0]=1.37726360474278; $bf[16,1]=0.1712060794402248; $bf[16,2]=0.1284154983726964; $bf[16,3]=0.1168828393481309; $bf[16,4]=1.554360452402591;
$bf[17,0]=1.370027884281683; $bf[17,1]=0.1700851860489365; $bf[17,2]=0.1283842567878963; $bf[17,3]=0.1169317269942427; $bf[17,4]=1.554410095172159;
$bf[18,0]=1.402603467659901; $bf[18,1]=0.1751256195758799; $bf[18,2]=0.128526939286985; $bf[18,3]=0.1167134270893575; $bf[18,4]=1.554238850961895;
$bf[19,0]=1.512291964877249; $bf[19,1]=0.1919951376210633; $bf[19,2]=0.1290428000295734; $bf[19,3]=0.1160096726882482; $bf[19,4]=1.554575771272056;
$bf[20,0]=1.584149839067948; $bf[20,1]=0.2029705700868192; $bf[20,2]=0.1294069535694362; $bf[20,3]=0.115571811609156; $bf[20,4]=1.555475393037952;
$bf[21,0]=1.545782890311557; $bf[21,1]=0.1971173836429614; $bf[21,2]=0.1292101376163686; $bf[21,3]=0.1158034892375213; $bf[21,4]=1.554933113257562;
$bf[22,0]=1.400813652664795; $bf[22,1]=0.1748490678692178; $bf[22,2]=0.1285189662205326; $bf[22,3]=0.1167253031910009; $bf[22,4]=1.554244821564239;
$bf[23,0]=1.245222641999817; $bf[23,1]=0.1506207781314117; $bf[23,2]=0.1278904877043274; $bf[23,3]=0.1178148008455686; $bf[23,4]=1.556422647054774;
$ws.Range("B2:F25").Value = $bf

# I2:J25
$ij = New-Object 'object[,]' 24,2
$ij[0,0]=0.9378677532000133; $ij[0,1]=0.1443678292441852;
$ij[1,0]=0.9455961383165139; $ij[1,1]=0.1451163965833979;
$ij[2,0]=0.9507960789339478; $ij[2,1]=0.1456016224407435;
$ij[3,0]=0.9530294468627538; $ij[3,1]=0.1458058067498973;
$ij[4,0]=0.9534072028578713; $ij[4,1]=0.1458401013926416;
$ij[5,0]=0.9508257358948917; $ij[5,1]=0.1456043500055655;
$ij[6,0]=0.940438141265453; $ij[6,1]=0.1446206292556642;
$ij[7,0]=0.9236752209440553; $ij[7,1]=0.1428941046335126;
$ij[8,0]=0.9135579864368992; $ij[8,1]=0.1417482617222534;
$ij[9,0]=0.9094327000760885; $ij[9,1]=0.1412534373732832;
$ij[10,0]=0.9079391637939125; $ij[10,1]=0.1410698470255136;
$ij[11,0]=0.9082577714238056; $ij[11,1]=0.1411092181543223;
$ij[12,0]=0.9093084506140201; $ij[12,1]=0.1412382574128674;
$ij[13,0]=0.9099609587556614; $ij[13,1]=0.1413177907845764;
$ij[14,0]=0.9138371852018992; $ij[14,1]=0.1417811304564207;
$ij[15,0]=0.9163373302082292; $ij[15,1]=0.1420721346129783;
$ij[16,0]=0.9178202529901753; $ij[16,1]=0.1422420004380989;
$ij[17,0]=0.9183300572871502; $ij[17,1]=0.1422999417304158;
$ij[18,0]=0.916066537847307; $ij[18,1]=0.1420408993126099;
$ij[19,0]=0.9089979782200572; $ij[19,1]=0.1412002527127134;
$ij[20,0]=0.9047782433379226; $ij[20,1]=0.1406729210891831;
$ij[21,0]=0.9069937924700326; $ij[21,1]=0.140952351242202;
$ij[22,0]=0.9161888211694986; $ij[22,1]=0.1420550128053648;
$ij[23,0]=0.9278239220136726; $ij[23,1]=0.1433395787458567;
$ws.Range("I2:J25").Value = $ij

# L2:O25
$lo = New-Object 'object[,]' 24,4
$lo[0,0]=0.3127432918536357; $lo[0,1]=0.2926087856766273; $lo[0,2]=1.574602123217531; $lo[0,3]=3.989614806516471;
$lo[1,0]=0.3096940352696365; $lo[1,1]=0.2799589163362128; $lo[1,2]=1.58738224020091; $lo[1,3]=4.002590084295832;
$lo[2,0]=0.3079292087611449; $lo[2,1]=0.2722774406975574; $lo[2,2]=1.595759612943993; $lo[2,3]=4.012665513329068;
$lo[3,0]=0.3072371620956105; $lo[3,1]=0.2691689674518187; $lo[3,2]=1.599306940158947; $lo[3,3]=4.01730171439462;
$lo[4,0]=0.3071238906754417; $lo[4,1]=0.268654131783812; $lo[4,2]=1.599904036709873; $lo[4,3]=4.018103590808408;
$lo[5,0]=0.3079197655514392; $lo[5,1]=0.2722354301577283; $lo[5,2]=1.595806912818446; $lo[5,3]=4.01272589104164;
$lo[6,0]=0.3116696819327132; $lo[6,1]=0.2882295010502887; $lo[6,2]=1.57889871205608; $lo[6,3]=3.993651145431784;
$lo[7,0]=0.3198709853739246; $lo[7,1]=0.3202629921164686; $lo[7,2]=1.549944137647238; $lo[7,3]=3.972972207443746;
$lo[8,0]=0.3264079081921665; $lo[8,1]=0.344195055063885; $lo[8,2]=1.531225424700601; $lo[8,3]=3.967976406767775;
$lo[9,0]=0.3294916942243162; $lo[9,1]=0.3551664361443088; $lo[9,2]=1.523262869450697; $lo[9,3]=3.967918270626711;
$lo[10,0]=0.3306751679384234; $lo[10,1]=0.3593329363468882; $lo[10,2]=1.520327028813675; $lo[10,3]=3.968214641128043;
$lo[11,0]=0.3304195889008241; $lo[11,1]=0.3584350822989251; $lo[11,2]=1.520955784061385; $lo[11,3]=3.968136653033611;
$lo[12,0]=0.3295887451501898; $lo[12,1]=0.3555089802311429; $lo[12,2]=1.523019745034674; $lo[12,3]=3.967936272602117;
$lo[13,0]=0.3290818717526349; $lo[13,1]=0.3537181964430047; $lo[13,2]=1.524294319791643; $lo[13,3]=3.96785499535693;
$lo[14,0]=0.3262085806917128; $lo[14,1]=0.343479726887324; $lo[14,2]=1.531756911810369; $lo[14,3]=3.968024762489762;
$lo[15,0]=0.3244740236079338; $lo[15,1]=0.3372202144390855; $lo[15,2]=1.536476475085706; $lo[15,3]=3.968696052248447;
$lo[16,0]=0.3234867233537813; $lo[16,1]=0.333627885113529; $lo[16,2]=1.539243067863836; $lo[16,3]=3.9692906008585;
$lo[17,0]=0.3231542254307271; $lo[17,1]=0.3324129627944572; $lo[17,2]=1.540188726112483; $lo[17,3]=3.969527706289824;
$lo[18,0]=0.3246575976877608; $lo[18,1]=0.3378857267281532; $lo[18,2]=1.535968685654197; $lo[18,3]=3.968603019881044;
$lo[19,0]=0.3298323586858913; $lo[19,1]=0.3563681272816055; $lo[19,2]=1.522411355130849; $lo[19,3]=3.967986488740024;
$lo[20,0]=0.3333058824518815; $lo[20,1]=0.3685165370189267; $lo[20,2]=1.514013667119727; $lo[20,3]=3.969439317836219;
$lo[21,0]=0.3314436639050768; $lo[21,1]=0.3620264739255177; $lo[21,2]=1.518453345653441; $lo[21,3]=3.968494130871676;
$lo[22,0]=0.3245745729385732; $lo[22,1]=0.3375848288279855; $lo[22,2]=1.536198091369634; $lo[22,3]=3.968644429984266;
$lo[23,0]=0.3175621165232769; $lo[23,1]=0.3115265719109246; $lo[23,2]=1.557328033575359; $lo[23,3]=3.976775706613552;
$ws.Range("L2:O25").Value = $lo
